$wb = $excel.ActiveWorkbook

# Rename the first sheet from "1SPL01_plants" to "plant_growth" (template name,
# dropping the "naming convention" prefix) per the commit message.
$wsPlants = $wb.Worksheets.Item(1)
$wsPlants.Name = "plant_growth"

# The other sheets keep their names.
$wsMetadata      = $wb.Worksheets.Item(2)   # SwateTemplateMetadata
$wsMetabolomics  = $wb.Worksheets.Item(3)   # METABOLIGHTS_METABOLOMICS

# Touch the metadata sheet so its stale scroll position ("topLeftCell=A13")
# is cleared, while its selection (F36) is unaffected.
$wsMetadata.Activate()

# Finally make the renamed sheet the active / selected tab, matching the
# updated workbookView (tabSelected moves from METABOLIGHTS_METABOLOMICS to
# plant_growth, and activeTab goes back to the first sheet).
$wsPlants.Activate()
